# Implementation of RQ4 for both regular and micro-clones are completed.

$wb = $excel.ActiveWorkbook

# --- RQ3 sheet: fill in the "All" row (row 3) with Regular-Clones numbers ---
$ws3 = $wb.Worksheets.Item("RQ3")
$ws3.Range("B3").Value = 16
$ws3.Range("C3").Value = 16
$ws3.Range("E3").Value = 129
$ws3.Range("F3").Value = 129

# --- RQ4 sheet: fill in numbers for rows 3 (All) and 4 (Ctags) ---
$ws4 = $wb.Worksheets.Item("RQ4")
$ws4.Range("B3").Value = 16
$ws4.Range("E3").Value = 129

$ws4.Range("B4").Value = 12
$ws4.Range("C4").Value = 2
$ws4.Range("E4").Value = 72
$ws4.Range("F4").Value = 10

# --- Update active tab / selection to match final state (RQ4 active, RQ3 not) ---
$ws3.Activate()
$ws3.Range("E3").Select()

$ws4.Activate()
$ws4.Range("F4").Select()
